$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-26 Saturday" "2024-10-27 Sunday"

Replace-Text "221×8=1768" "401×7=2807"
Replace-Text "455×2=910" "552×7=3864"
Replace-Text "295×2=590" "192×7=1344"
Replace-Text "818×7=5726" "568×5=2840"
Replace-Text "350×4=1400" "348×3=1044"

Replace-Text "969×3=2907" "564×8=4512"
Replace-Text "979×9=8811" "898×4=3592"
Replace-Text "287×7=2009" "627×4=2508"
Replace-Text "770×3=2310" "290×6=1740"
Replace-Text "195×7=1365" "567×9=5103"

Replace-Text "608×4=2432" "540×3=1620"
Replace-Text "649×5=3245" "298×9=2682"
Replace-Text "207×9=1863" "674×4=2696"
Replace-Text "924×9=8316" "762×9=6858"
Replace-Text "624×9=5616" "279×3=837"

Replace-Text "876×5=4380" "705×2=1410"
Replace-Text "720×7=5040" "606×9=5454"
Replace-Text "881×5=4405" "578×3=1734"
Replace-Text "161×9=1449" "387×9=3483"
Replace-Text "830×5=4150" "294×4=1176"

Replace-Text "354×2=708" "637×9=5733"
Replace-Text "321×8=2568" "206×4=824"
Replace-Text "741×2=1482" "251×3=753"
Replace-Text "437×3=1311" "595×5=2975"
Replace-Text "349×6=2094" "558×2=1116"
